# Apply policy update to the Vermont policy effectiveness workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (rows 18-101): policy weight updated from 0.5 to 1
$ws.Range("B18:B101").Value = 1

# Column H (rows 26-59): policy flag turned off (1 -> 0)
$ws.Range("H26:H59").Value = 0
